# Fill A1:A5 with the numbers 1-5 (matches the new <sheetData> rows in the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 5; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# After typing, Excel leaves the selection one row below the last entry (A6),
# matching the <selection activeCell="A6" sqref="A6"/> in the target sheetView.
$ws.Range("A6").Select()

# Print setup: paperSize="9" (A4) / orientation="portrait" on the sheet's <pageSetup>.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
